# Auto-generated Excel COM-interop script
# Applies the "river update May 2024" edit to the Tiraumea at Ngaturi sheet:
#  1) Updates 70 recomputed statistic cells (Mean/G column & a few summer-period values)
#     in the existing 2009-2013 .. 2018-2022 year-range blocks (rows 3-197).
#  2) Appends a new 20-row "2019 - 2023" year-range block (rows 207-226).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: update recomputed values in existing rows ---
$ws.Cells.Item(3,7).Value = 0.820590703987405
$ws.Cells.Item(4,7).Value = 0.0111808623123795
$ws.Cells.Item(5,7).Value = 0.0111808623123795
$ws.Cells.Item(19,7).Value = 0.0646378193482504
$ws.Cells.Item(20,7).Value = 0.0646378193482504
$ws.Cells.Item(24,7).Value = 0.0104791938456179
$ws.Cells.Item(25,7).Value = 0.0104791938456179
$ws.Cells.Item(40,7).Value = 0.0644165587033087
$ws.Cells.Item(41,7).Value = 0.0644165587033087
$ws.Cells.Item(45,7).Value = 0.0095849325406439
$ws.Cells.Item(46,7).Value = 0.0095849325406439
$ws.Cells.Item(52,7).Value = 0.0152318746332107
$ws.Cells.Item(52,12).Value = 0.00257
$ws.Cells.Item(53,7).Value = 0.0152318746332107
$ws.Cells.Item(53,12).Value = 0.00257
$ws.Cells.Item(61,7).Value = 0.0522998920366421
$ws.Cells.Item(62,7).Value = 0.0522998920366421
$ws.Cells.Item(66,7).Value = 0.0101181110423397
$ws.Cells.Item(67,7).Value = 0.0101181110423397
$ws.Cells.Item(73,6).Value = 0.00661
$ws.Cells.Item(73,7).Value = 0.0122995037553796
$ws.Cells.Item(73,12).Value = 0.00257
$ws.Cells.Item(74,6).Value = 0.00661
$ws.Cells.Item(74,7).Value = 0.0122995037553796
$ws.Cells.Item(74,12).Value = 0.00257
$ws.Cells.Item(87,7).Value = 0.0111672199334733
$ws.Cells.Item(88,7).Value = 0.0111672199334733
$ws.Cells.Item(94,6).Value = 0.00607
$ws.Cells.Item(94,7).Value = 0.0117451718234435
$ws.Cells.Item(94,12).Value = 0.00257
$ws.Cells.Item(95,6).Value = 0.00607
$ws.Cells.Item(95,7).Value = 0.0117451718234435
$ws.Cells.Item(95,12).Value = 0.00257
$ws.Cells.Item(108,7).Value = 0.0123005532668067
$ws.Cells.Item(109,7).Value = 0.0123005532668067
$ws.Cells.Item(115,7).Value = 0.0144496875245774
$ws.Cells.Item(115,12).Value = 0.00591
$ws.Cells.Item(116,7).Value = 0.0144496875245774
$ws.Cells.Item(116,12).Value = 0.00591
$ws.Cells.Item(128,7).Value = 0.799897674280284
$ws.Cells.Item(128,9).Value = 2.30758
$ws.Cells.Item(128,14).Value = 2.1754
$ws.Cells.Item(129,7).Value = 0.0128657011755503
$ws.Cells.Item(130,7).Value = 0.0128657011755503
$ws.Cells.Item(136,7).Value = 0.0150008566457504
$ws.Cells.Item(137,7).Value = 0.0150008566457504
$ws.Cells.Item(148,7).Value = 0.839365614484559
$ws.Cells.Item(148,9).Value = 2.2
$ws.Cells.Item(148,14).Value = 2.1706
$ws.Cells.Item(149,7).Value = 0.0133937711902705
$ws.Cells.Item(150,7).Value = 0.0133937711902705
$ws.Cells.Item(156,7).Value = 0.0177721683358035
$ws.Cells.Item(157,7).Value = 0.0177721683358035
$ws.Cells.Item(168,7).Value = 0.7625512851395549
$ws.Cells.Item(168,8).Value = 2.33447441113531
$ws.Cells.Item(168,9).Value = 2.17
$ws.Cells.Item(168,13).Value = 1.8
$ws.Cells.Item(168,14).Value = 1.926
$ws.Cells.Item(169,7).Value = 0.0127570495607693
$ws.Cells.Item(170,7).Value = 0.0127570495607693
$ws.Cells.Item(176,7).Value = 0.0181046490819467
$ws.Cells.Item(177,7).Value = 0.0181046490819467
$ws.Cells.Item(188,7).Value = 0.7516550497641999
$ws.Cells.Item(188,8).Value = 2.33447441113531
$ws.Cells.Item(188,9).Value = 2.2
$ws.Cells.Item(188,14).Value = 1.958
$ws.Cells.Item(189,7).Value = 0.0131074340822658
$ws.Cells.Item(190,7).Value = 0.0131074340822658
$ws.Cells.Item(196,7).Value = 0.0181933055413241
$ws.Cells.Item(197,7).Value = 0.0181933055413241

# --- Part 2: append new "2019 - 2023" block (rows 207-226) ---
# row 207
$ws.Cells.Item(207,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(207,2).Value = "ASPM"
$ws.Cells.Item(207,3).Value = "C"
$ws.Cells.Item(207,4).Value = "2019 - 2023"
$ws.Cells.Item(207,5).Value = "RepSite"
$ws.Cells.Item(207,6).Value = 0.392
$ws.Cells.Item(207,7).Value = 0.3762
$ws.Cells.Item(207,8).Value = 0.41
$ws.Cells.Item(207,9).Value = 0.41
$ws.Cells.Item(207,12).Value = 0.3765
$ws.Cells.Item(207,13).Value = 0.4065
$ws.Cells.Item(207,14).Value = 0.41
$ws.Cells.Item(207,15).Value = 1847797
$ws.Cells.Item(207,16).Value = 5516292
$ws.Cells.Item(207,17).Value = "Tararua District"
$ws.Cells.Item(207,18).Value = "Manawatū"
$ws.Cells.Item(207,19).Value = "Tiraumea"
$ws.Cells.Item(207,20).Value = "Mana_7b"

# row 208
$ws.Cells.Item(208,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(208,2).Value = "Visual Clarity (Sediment class 4)"
$ws.Cells.Item(208,3).Value = "D"
$ws.Cells.Item(208,4).Value = "2019 - 2023"
$ws.Cells.Item(208,5).Value = "RepSite"
$ws.Cells.Item(208,6).Value = 0.3
$ws.Cells.Item(208,7).Value = 0.685023857357517
$ws.Cells.Item(208,8).Value = 2.33447441113531
$ws.Cells.Item(208,9).Value = 2.2
$ws.Cells.Item(208,12).Value = 1.15
$ws.Cells.Item(208,13).Value = 1.8
$ws.Cells.Item(208,14).Value = 1.982
$ws.Cells.Item(208,15).Value = 1847797
$ws.Cells.Item(208,16).Value = 5516292
$ws.Cells.Item(208,17).Value = "Tararua District"
$ws.Cells.Item(208,18).Value = "Manawatū"
$ws.Cells.Item(208,19).Value = "Tiraumea"
$ws.Cells.Item(208,20).Value = "Mana_7b"
$ws.Cells.Item(208,21).Value = "m"

# row 209
$ws.Cells.Item(209,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(209,2).Value = "DRP (95th Percentile)"
$ws.Cells.Item(209,3).Value = "B"
$ws.Cells.Item(209,4).Value = "2019 - 2023"
$ws.Cells.Item(209,5).Value = "RepSite"
$ws.Cells.Item(209,6).Value = 0.012
$ws.Cells.Item(209,7).Value = 0.0125074340822658
$ws.Cells.Item(209,8).Value = 0.026
$ws.Cells.Item(209,9).Value = 0.024
$ws.Cells.Item(209,12).Value = 0.008999999999999999
$ws.Cells.Item(209,13).Value = 0.019
$ws.Cells.Item(209,14).Value = 0.0231
$ws.Cells.Item(209,15).Value = 1847797
$ws.Cells.Item(209,16).Value = 5516292
$ws.Cells.Item(209,17).Value = "Tararua District"
$ws.Cells.Item(209,18).Value = "Manawatū"
$ws.Cells.Item(209,19).Value = "Tiraumea"
$ws.Cells.Item(209,20).Value = "Mana_7b"
$ws.Cells.Item(209,21).Value = "mg/L"

# row 210
$ws.Cells.Item(210,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(210,2).Value = "DRP (Median)"
$ws.Cells.Item(210,3).Value = "C"
$ws.Cells.Item(210,4).Value = "2019 - 2023"
$ws.Cells.Item(210,5).Value = "RepSite"
$ws.Cells.Item(210,6).Value = 0.012
$ws.Cells.Item(210,7).Value = 0.0125074340822658
$ws.Cells.Item(210,8).Value = 0.026
$ws.Cells.Item(210,9).Value = 0.024
$ws.Cells.Item(210,12).Value = 0.008999999999999999
$ws.Cells.Item(210,13).Value = 0.019
$ws.Cells.Item(210,14).Value = 0.0231
$ws.Cells.Item(210,15).Value = 1847797
$ws.Cells.Item(210,16).Value = 5516292
$ws.Cells.Item(210,17).Value = "Tararua District"
$ws.Cells.Item(210,18).Value = "Manawatū"
$ws.Cells.Item(210,19).Value = "Tiraumea"
$ws.Cells.Item(210,20).Value = "Mana_7b"
$ws.Cells.Item(210,21).Value = "mg/L"

# row 211
$ws.Cells.Item(211,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(211,2).Value = "E coli (>260)"
$ws.Cells.Item(211,3).Value = "E"
$ws.Cells.Item(211,4).Value = "2019 - 2023"
$ws.Cells.Item(211,5).Value = "RepSite"
$ws.Cells.Item(211,6).Value = 380
$ws.Cells.Item(211,7).Value = 3670.33333333333
$ws.Cells.Item(211,8).Value = 98000
$ws.Cells.Item(211,9).Value = 14472
$ws.Cells.Item(211,10).Value = 46.2962962962963
$ws.Cells.Item(211,11).Value = 59.2592592592593
$ws.Cells.Item(211,12).Value = 185
$ws.Cells.Item(211,13).Value = 3300
$ws.Cells.Item(211,14).Value = 8702.799999999999
$ws.Cells.Item(211,15).Value = 1847797
$ws.Cells.Item(211,16).Value = 5516292
$ws.Cells.Item(211,17).Value = "Tararua District"
$ws.Cells.Item(211,18).Value = "Manawatū"
$ws.Cells.Item(211,19).Value = "Tiraumea"
$ws.Cells.Item(211,20).Value = "Mana_7b"
$ws.Cells.Item(211,21).Value = "% exceedances over 260/100 mL"

# row 212
$ws.Cells.Item(212,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(212,2).Value = "E coli (>540)"
$ws.Cells.Item(212,3).Value = "E"
$ws.Cells.Item(212,4).Value = "2019 - 2023"
$ws.Cells.Item(212,5).Value = "RepSite"
$ws.Cells.Item(212,6).Value = 380
$ws.Cells.Item(212,7).Value = 3670.33333333333
$ws.Cells.Item(212,8).Value = 98000
$ws.Cells.Item(212,9).Value = 14472
$ws.Cells.Item(212,10).Value = 46.2962962962963
$ws.Cells.Item(212,11).Value = 59.2592592592593
$ws.Cells.Item(212,12).Value = 185
$ws.Cells.Item(212,13).Value = 3300
$ws.Cells.Item(212,14).Value = 8702.799999999999
$ws.Cells.Item(212,15).Value = 1847797
$ws.Cells.Item(212,16).Value = 5516292
$ws.Cells.Item(212,17).Value = "Tararua District"
$ws.Cells.Item(212,18).Value = "Manawatū"
$ws.Cells.Item(212,19).Value = "Tiraumea"
$ws.Cells.Item(212,20).Value = "Mana_7b"
$ws.Cells.Item(212,21).Value = "% exceedances over 540/100 mL"

# row 213
$ws.Cells.Item(213,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(213,2).Value = "E coli (Median)"
$ws.Cells.Item(213,3).Value = "E"
$ws.Cells.Item(213,4).Value = "2019 - 2023"
$ws.Cells.Item(213,5).Value = "RepSite"
$ws.Cells.Item(213,6).Value = 380
$ws.Cells.Item(213,7).Value = 3670.33333333333
$ws.Cells.Item(213,8).Value = 98000
$ws.Cells.Item(213,9).Value = 14472
$ws.Cells.Item(213,10).Value = 46.2962962962963
$ws.Cells.Item(213,11).Value = 59.2592592592593
$ws.Cells.Item(213,12).Value = 185
$ws.Cells.Item(213,13).Value = 3300
$ws.Cells.Item(213,14).Value = 8702.799999999999
$ws.Cells.Item(213,15).Value = 1847797
$ws.Cells.Item(213,16).Value = 5516292
$ws.Cells.Item(213,17).Value = "Tararua District"
$ws.Cells.Item(213,18).Value = "Manawatū"
$ws.Cells.Item(213,19).Value = "Tiraumea"
$ws.Cells.Item(213,20).Value = "Mana_7b"
$ws.Cells.Item(213,21).Value = "E. coli/100 mL"

# row 214
$ws.Cells.Item(214,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(214,2).Value = "E coli (95th Percentile)"
$ws.Cells.Item(214,3).Value = "E"
$ws.Cells.Item(214,4).Value = "2019 - 2023"
$ws.Cells.Item(214,5).Value = "RepSite"
$ws.Cells.Item(214,6).Value = 380
$ws.Cells.Item(214,7).Value = 3670.33333333333
$ws.Cells.Item(214,8).Value = 98000
$ws.Cells.Item(214,9).Value = 14472
$ws.Cells.Item(214,10).Value = 46.2962962962963
$ws.Cells.Item(214,11).Value = 59.2592592592593
$ws.Cells.Item(214,12).Value = 185
$ws.Cells.Item(214,13).Value = 3300
$ws.Cells.Item(214,14).Value = 8702.799999999999
$ws.Cells.Item(214,15).Value = 1847797
$ws.Cells.Item(214,16).Value = 5516292
$ws.Cells.Item(214,17).Value = "Tararua District"
$ws.Cells.Item(214,18).Value = "Manawatū"
$ws.Cells.Item(214,19).Value = "Tiraumea"
$ws.Cells.Item(214,20).Value = "Mana_7b"
$ws.Cells.Item(214,21).Value = "E. coli/100 mL"

# row 215
$ws.Cells.Item(215,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(215,2).Value = "MCI"
$ws.Cells.Item(215,3).Value = "C"
$ws.Cells.Item(215,4).Value = "2019 - 2023"
$ws.Cells.Item(215,5).Value = "RepSite"
$ws.Cells.Item(215,6).Value = 106.32
$ws.Cells.Item(215,7).Value = 107.198
$ws.Cells.Item(215,8).Value = 113
$ws.Cells.Item(215,9).Value = 113
$ws.Cells.Item(215,12).Value = 108.335
$ws.Cells.Item(215,13).Value = 112.5345
$ws.Cells.Item(215,14).Value = 113
$ws.Cells.Item(215,15).Value = 1847797
$ws.Cells.Item(215,16).Value = 5516292
$ws.Cells.Item(215,17).Value = "Tararua District"
$ws.Cells.Item(215,18).Value = "Manawatū"
$ws.Cells.Item(215,19).Value = "Tiraumea"
$ws.Cells.Item(215,20).Value = "Mana_7b"

# row 216
$ws.Cells.Item(216,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(216,2).Value = "Ammoniacal-N (95th Percentile)"
$ws.Cells.Item(216,3).Value = "A"
$ws.Cells.Item(216,4).Value = "2019 - 2023"
$ws.Cells.Item(216,5).Value = "RepSite"
$ws.Cells.Item(216,6).Value = 0.0137
$ws.Cells.Item(216,7).Value = 0.016931562418892
$ws.Cells.Item(216,8).Value = 0.0628709597299737
$ws.Cells.Item(216,9).Value = 0.0433
$ws.Cells.Item(216,12).Value = 0.01457
$ws.Cells.Item(216,13).Value = 0.02936
$ws.Cells.Item(216,14).Value = 0.03823
$ws.Cells.Item(216,15).Value = 1847797
$ws.Cells.Item(216,16).Value = 5516292
$ws.Cells.Item(216,17).Value = "Tararua District"
$ws.Cells.Item(216,18).Value = "Manawatū"
$ws.Cells.Item(216,19).Value = "Tiraumea"
$ws.Cells.Item(216,20).Value = "Mana_7b"
$ws.Cells.Item(216,21).Value = "mg NH4-N/L"

# row 217
$ws.Cells.Item(217,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(217,2).Value = "Ammoniacal-N (Median)"
$ws.Cells.Item(217,3).Value = "A"
$ws.Cells.Item(217,4).Value = "2019 - 2023"
$ws.Cells.Item(217,5).Value = "RepSite"
$ws.Cells.Item(217,6).Value = 0.0137
$ws.Cells.Item(217,7).Value = 0.016931562418892
$ws.Cells.Item(217,8).Value = 0.0628709597299737
$ws.Cells.Item(217,9).Value = 0.0433
$ws.Cells.Item(217,12).Value = 0.01457
$ws.Cells.Item(217,13).Value = 0.02936
$ws.Cells.Item(217,14).Value = 0.03823
$ws.Cells.Item(217,15).Value = 1847797
$ws.Cells.Item(217,16).Value = 5516292
$ws.Cells.Item(217,17).Value = "Tararua District"
$ws.Cells.Item(217,18).Value = "Manawatū"
$ws.Cells.Item(217,19).Value = "Tiraumea"
$ws.Cells.Item(217,20).Value = "Mana_7b"
$ws.Cells.Item(217,21).Value = "mg NH4-N/L"

# row 218
$ws.Cells.Item(218,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(218,2).Value = "Nitrate-N (95th Percentile)"
$ws.Cells.Item(218,3).Value = "A"
$ws.Cells.Item(218,4).Value = "2019 - 2023"
$ws.Cells.Item(218,5).Value = "RepSite"
$ws.Cells.Item(218,6).Value = 0.615
$ws.Cells.Item(218,7).Value = 0.628254545454545
$ws.Cells.Item(218,8).Value = 1.06
$ws.Cells.Item(218,9).Value = 0.956
$ws.Cells.Item(218,12).Value = 0.6
$ws.Cells.Item(218,13).Value = 0.72995
$ws.Cells.Item(218,14).Value = 0.8732
$ws.Cells.Item(218,15).Value = 1847797
$ws.Cells.Item(218,16).Value = 5516292
$ws.Cells.Item(218,17).Value = "Tararua District"
$ws.Cells.Item(218,18).Value = "Manawatū"
$ws.Cells.Item(218,19).Value = "Tiraumea"
$ws.Cells.Item(218,20).Value = "Mana_7b"
$ws.Cells.Item(218,21).Value = "mg NO3-N/L"

# row 219
$ws.Cells.Item(219,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(219,2).Value = "Nitrate-N (Median)"
$ws.Cells.Item(219,3).Value = "A"
$ws.Cells.Item(219,4).Value = "2019 - 2023"
$ws.Cells.Item(219,5).Value = "RepSite"
$ws.Cells.Item(219,6).Value = 0.615
$ws.Cells.Item(219,7).Value = 0.628254545454545
$ws.Cells.Item(219,8).Value = 1.06
$ws.Cells.Item(219,9).Value = 0.956
$ws.Cells.Item(219,12).Value = 0.6
$ws.Cells.Item(219,13).Value = 0.72995
$ws.Cells.Item(219,14).Value = 0.8732
$ws.Cells.Item(219,15).Value = 1847797
$ws.Cells.Item(219,16).Value = 5516292
$ws.Cells.Item(219,17).Value = "Tararua District"
$ws.Cells.Item(219,18).Value = "Manawatū"
$ws.Cells.Item(219,19).Value = "Tiraumea"
$ws.Cells.Item(219,20).Value = "Mana_7b"
$ws.Cells.Item(219,21).Value = "mg NO3-N/L"

# row 220
$ws.Cells.Item(220,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(220,2).Value = "QMCI"
$ws.Cells.Item(220,3).Value = "D"
$ws.Cells.Item(220,4).Value = "2019 - 2023"
$ws.Cells.Item(220,5).Value = "RepSite"
$ws.Cells.Item(220,6).Value = 4.29
$ws.Cells.Item(220,7).Value = 4.501
$ws.Cells.Item(220,8).Value = 5.554
$ws.Cells.Item(220,9).Value = 5.554
$ws.Cells.Item(220,12).Value = 4.07
$ws.Cells.Item(220,13).Value = 5.3531
$ws.Cells.Item(220,14).Value = 5.554
$ws.Cells.Item(220,15).Value = 1847797
$ws.Cells.Item(220,16).Value = 5516292
$ws.Cells.Item(220,17).Value = "Tararua District"
$ws.Cells.Item(220,18).Value = "Manawatū"
$ws.Cells.Item(220,19).Value = "Tiraumea"
$ws.Cells.Item(220,20).Value = "Mana_7b"

# row 221
$ws.Cells.Item(221,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(221,2).Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Cells.Item(221,4).Value = "2019 - 2023"
$ws.Cells.Item(221,5).Value = "RepSite"
$ws.Cells.Item(221,6).Value = 0.65
$ws.Cells.Item(221,7).Value = 0.658509090909091
$ws.Cells.Item(221,8).Value = 1.112
$ws.Cells.Item(221,9).Value = 1.01325
$ws.Cells.Item(221,12).Value = 0.642
$ws.Cells.Item(221,13).Value = 0.74885
$ws.Cells.Item(221,14).Value = 0.9042
$ws.Cells.Item(221,15).Value = 1847797
$ws.Cells.Item(221,16).Value = 5516292
$ws.Cells.Item(221,17).Value = "Tararua District"
$ws.Cells.Item(221,18).Value = "Manawatū"
$ws.Cells.Item(221,19).Value = "Tiraumea"
$ws.Cells.Item(221,20).Value = "Mana_7b"
$ws.Cells.Item(221,21).Value = "g/m3"

# row 222
$ws.Cells.Item(222,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(222,2).Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Cells.Item(222,4).Value = "2019 - 2023"
$ws.Cells.Item(222,5).Value = "RepSite"
$ws.Cells.Item(222,6).Value = 0.65
$ws.Cells.Item(222,7).Value = 0.658509090909091
$ws.Cells.Item(222,8).Value = 1.112
$ws.Cells.Item(222,9).Value = 1.01325
$ws.Cells.Item(222,12).Value = 0.642
$ws.Cells.Item(222,13).Value = 0.74885
$ws.Cells.Item(222,14).Value = 0.9042
$ws.Cells.Item(222,15).Value = 1847797
$ws.Cells.Item(222,16).Value = 5516292
$ws.Cells.Item(222,17).Value = "Tararua District"
$ws.Cells.Item(222,18).Value = "Manawatū"
$ws.Cells.Item(222,19).Value = "Tiraumea"
$ws.Cells.Item(222,20).Value = "Mana_7b"
$ws.Cells.Item(222,21).Value = "g/m3"

# row 223
$ws.Cells.Item(223,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(223,2).Value = "Total Nitrogen (95th Percentile)"
$ws.Cells.Item(223,4).Value = "2019 - 2023"
$ws.Cells.Item(223,5).Value = "RepSite"
$ws.Cells.Item(223,6).Value = 0.9
$ws.Cells.Item(223,7).Value = 1.08109090909091
$ws.Cells.Item(223,8).Value = 3.03
$ws.Cells.Item(223,9).Value = 2.1275
$ws.Cells.Item(223,12).Value = 0.85
$ws.Cells.Item(223,13).Value = 1.473
$ws.Cells.Item(223,14).Value = 1.859
$ws.Cells.Item(223,15).Value = 1847797
$ws.Cells.Item(223,16).Value = 5516292
$ws.Cells.Item(223,17).Value = "Tararua District"
$ws.Cells.Item(223,18).Value = "Manawatū"
$ws.Cells.Item(223,19).Value = "Tiraumea"
$ws.Cells.Item(223,20).Value = "Mana_7b"
$ws.Cells.Item(223,21).Value = "g/m3"

# row 224
$ws.Cells.Item(224,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(224,2).Value = "Total Nitrogen (Median)"
$ws.Cells.Item(224,4).Value = "2019 - 2023"
$ws.Cells.Item(224,5).Value = "RepSite"
$ws.Cells.Item(224,6).Value = 0.9
$ws.Cells.Item(224,7).Value = 1.08109090909091
$ws.Cells.Item(224,8).Value = 3.03
$ws.Cells.Item(224,9).Value = 2.1275
$ws.Cells.Item(224,12).Value = 0.85
$ws.Cells.Item(224,13).Value = 1.473
$ws.Cells.Item(224,14).Value = 1.859
$ws.Cells.Item(224,15).Value = 1847797
$ws.Cells.Item(224,16).Value = 5516292
$ws.Cells.Item(224,17).Value = "Tararua District"
$ws.Cells.Item(224,18).Value = "Manawatū"
$ws.Cells.Item(224,19).Value = "Tiraumea"
$ws.Cells.Item(224,20).Value = "Mana_7b"
$ws.Cells.Item(224,21).Value = "g/m3"

# row 225
$ws.Cells.Item(225,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(225,2).Value = "Total Phosphorus (95th Percentile)"
$ws.Cells.Item(225,4).Value = "2019 - 2023"
$ws.Cells.Item(225,5).Value = "RepSite"
$ws.Cells.Item(225,6).Value = 0.048
$ws.Cells.Item(225,7).Value = 0.1406
$ws.Cells.Item(225,8).Value = 1.74
$ws.Cells.Item(225,9).Value = 0.67525
$ws.Cells.Item(225,12).Value = 0.018
$ws.Cells.Item(225,13).Value = 0.2109
$ws.Cells.Item(225,14).Value = 0.3586
$ws.Cells.Item(225,15).Value = 1847797
$ws.Cells.Item(225,16).Value = 5516292
$ws.Cells.Item(225,17).Value = "Tararua District"
$ws.Cells.Item(225,18).Value = "Manawatū"
$ws.Cells.Item(225,19).Value = "Tiraumea"
$ws.Cells.Item(225,20).Value = "Mana_7b"
$ws.Cells.Item(225,21).Value = "g/m3"

# row 226
$ws.Cells.Item(226,1).Value = "Tiraumea at Ngaturi"
$ws.Cells.Item(226,2).Value = "Total Phosphorus (Median)"
$ws.Cells.Item(226,4).Value = "2019 - 2023"
$ws.Cells.Item(226,5).Value = "RepSite"
$ws.Cells.Item(226,6).Value = 0.048
$ws.Cells.Item(226,7).Value = 0.1406
$ws.Cells.Item(226,8).Value = 1.74
$ws.Cells.Item(226,9).Value = 0.67525
$ws.Cells.Item(226,12).Value = 0.018
$ws.Cells.Item(226,13).Value = 0.2109
$ws.Cells.Item(226,14).Value = 0.3586
$ws.Cells.Item(226,15).Value = 1847797
$ws.Cells.Item(226,16).Value = 5516292
$ws.Cells.Item(226,17).Value = "Tararua District"
$ws.Cells.Item(226,18).Value = "Manawatū"
$ws.Cells.Item(226,19).Value = "Tiraumea"
$ws.Cells.Item(226,20).Value = "Mana_7b"
$ws.Cells.Item(226,21).Value = "g/m3"

